$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.694.12'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.583.58'
$ws.Range('E3').Value = '  +2.06%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.17'
$ws.Range('E5').Value = '  +3.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.55'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.576.88'
$ws.Range('E7').Value = '  +2.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.617'
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  +6.34%  '
$ws.Range('E11').Value = '  +10.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.591'
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.22'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000279'
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.163.93'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.46'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '617.83'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.571.82'
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.834.03'
$ws.Range('E19').Value = '  +2.62%  '
$ws.Range('E20').Value = '  -2.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.50'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.891'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.32'
$ws.Range('E23').Value = '  -16.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.13'
$ws.Range('E24').Value = '  +1.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.90'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.81'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.66'
$ws.Range('E28').Value = '  +1.41%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.33'
$ws.Range('E29').Value = '  +0.74%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.93'
$ws.Range('E30').Value = '  +3.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.53'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.24'
$ws.Range('E33').Value = '  +4.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.31'
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '629.03'
$ws.Range('E35').Value = '  -0.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.79'
$ws.Range('E36').Value = '  +8.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.102'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.90'
$ws.Range('E38').Value = '  +1.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0487'
$ws.Range('E39').Value = '  +7.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '57.47'
$ws.Range('E40').Value = '  +0.43%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.405.01'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.325'
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.03'
$ws.Range('E45').Value = '  +10.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₃0719'
$ws.Range('E46').Value = '  +3.07%  '
$ws.Range('E47').Value = '  +7.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '33.07'
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.65'
$ws.Range('E50').Value = '  +0.36%  '
